$wb = $excel.ActiveWorkbook

$wsSite = $wb.Worksheets.Item("site_metrics")
$wsSite.Range("O3").Value = 0.006638589657142382
$wsSite.Range("O4").Value = 0.01689687697811004
$wsSite.Range("O11").Value = 0.3414642721151773
$wsSite.Range("O13").Value = 0.00446236154459127
$wsSite.Range("AK13").Value = $true
$wsSite.Range("O14").Value = 0.008533913737589298
$wsSite.Range("AK14").Value = $true
$wsSite.Range("O20").Value = 0.00396137615532286
$wsSite.Range("O21").Value = 0.09449115251611298
$wsSite.Range("O34").Value = 0.01688877555810498
$wsSite.Range("O36").Value = 0.0009109086966998175
$wsSite.Range("AK36").Value = $true
$wsSite.Range("AK37").Value = $true
$wsSite.Range("O39").Value = 0.353867549199511
$wsSite.Range("AK40").Value = $true
$wsSite.Range("AK41").Value = $true
$wsSite.Range("O42").Value = 0.5536757223203924
$wsSite.Range("O43").Value = 0.003720445005879518
$wsSite.Range("AK43").Value = $true
$wsSite.Range("AK44").Value = $true
$wsSite.Range("O52").Value = 0.0007259547207608328
$wsSite.Range("AK52").Value = $true
$wsSite.Range("O55").Value = 0.06873677721530226
$wsSite.Range("O56").Value = 0.08879613705747835
$wsSite.Range("O58").Value = 0.7306720040936679
$wsSite.Range("AK60").Value = $true
$wsSite.Range("O62").Value = 0.0009916035592921026
$wsSite.Range("O67").Value = 0.01099562600072403
$wsSite.Range("AK69").Value = $true
$wsSite.Range("O72").Value = 0.007682911412505886
$wsSite.Range("AK73").Value = $true
$wsSite.Range("O76").Value = 0.0009276591751838549
$wsSite.Range("O79").Value = 0.01676367222418533
$wsDur = $wb.Worksheets.Item("mk_duration")
$wsDur.Range("K5").Value = "no trend"
$wsDur.Range("L5").Value = $false
$wsDur.Range("M5").Value = 0.1544237061170961
$wsDur.Range("N5").Value = 1.424078649513432
$wsDur.Range("O5").Value = 0.3888888888888889
$wsDur.Range("P5").Value = 14
$wsDur.Range("Q5").Value = 83.33333333333333
$wsDur.Range("R5").Value = 0.9166666666666667
$wsDur.Range("S5").Value = -2.666666666666667
$wsDur.Range("M31").Value = 0.9311749148236654
$wsDur.Range("N31").Value = -0.08636670341750609
$wsDur.Range("O31").Value = -0.01231527093596059
$wsDur.Range("P31").Value = -5
$wsDur.Range("Q31").Value = 2145
$wsDur.Range("K61").Value = "no trend"
$wsDur.Range("L61").Value = $false
$wsDur.Range("M61").Value = 0.213590891281481
$wsDur.Range("N61").Value = 1.243751891458663
$wsDur.Range("O61").Value = 0.1897233201581028
$wsDur.Range("P61").Value = 48
$wsDur.Range("Q61").Value = 1428
$wsDur.Range("R61").Value = 0.08333333333333337
$wsDur.Range("S61").Value = 1.75
$wsDur.Range("M78").Value = 0.5807172923623485
$wsDur.Range("N78").Value = -0.5523372814706976
$wsDur.Range("O78").Value = -0.08
$wsDur.Range("P78").Value = -26
$wsDur.Range("Q78").Value = 2048.666666666667
$wsDur.Range("R78").Value = -0.02083333333333333
$wsDur.Range("S78").Value = 3.59375
$wsDur.Range("K82").Value = "no trend"
$wsDur.Range("L82").Value = $false
$wsDur.Range("M82").Value = 0.3513050661832442
$wsDur.Range("N82").Value = -0.9320608751309648
$wsDur.Range("O82").Value = -0.1978021978021978
$wsDur.Range("P82").Value = -18
$wsDur.Range("Q82").Value = 332.6666666666667
$wsDur.Range("R82").Value = -0.11875
$wsDur.Range("S82").Value = 3.549652777777778
$wsIntra = $wb.Worksheets.Item("mk_intra_annual")
$wsIntra.Range("M5").Value = 0.5084542305885602
$wsIntra.Range("N5").Value = 0.6612465225335806
$wsIntra.Range("O5").Value = 0.1944444444444444
$wsIntra.Range("P5").Value = 7
$wsIntra.Range("Q5").Value = 82.33333333333333
$wsIntra.Range("R5").Value = 0.08333333333333333
$wsIntra.Range("S5").Value = 0.6666666666666667
$wsIntra.Range("M31").Value = 0.9308309497786535
$wsIntra.Range("N31").Value = 0.08679941859050837
$wsIntra.Range("O31").Value = 0.01231527093596059
$wsIntra.Range("P31").Value = 5
$wsIntra.Range("Q31").Value = 2123.666666666667
$wsIntra.Range("M61").Value = 0.5742632659400388
$wsIntra.Range("N61").Value = -0.5617838953175571
$wsIntra.Range("O61").Value = -0.08695652173913043
$wsIntra.Range("P61").Value = -22
$wsIntra.Range("Q61").Value = 1397.333333333333
$wsIntra.Range("S61").Value = 5
$wsIntra.Range("M78").Value = 1
$wsIntra.Range("N78").Value = 0
$wsIntra.Range("O78").Value = 0.003076923076923077
$wsIntra.Range("P78").Value = 1
$wsIntra.Range("Q78").Value = 2015
$wsIntra.Range("R78").Value = 0
$wsIntra.Range("S78").Value = 4
$wsIntra.Range("K82").Value = "no trend"
$wsIntra.Range("L82").Value = $false
$wsIntra.Range("M82").Value = 0.866677633709956
$wsIntra.Range("N82").Value = 0.1678800645554932
$wsIntra.Range("O82").Value = 0.04395604395604396
$wsIntra.Range("P82").Value = 4
$wsIntra.Range("Q82").Value = 319.3333333333333
$wsIntra.Range("R82").Value = 0
$wsIntra.Range("S82").Value = 5
